# "Generate Report for handoff"
#
# A handoff attempt for e8f4dd36-ca04-4291-ab8e-ad32c76beac3.md failed during
# transform, so the localization-status report now reflects the NEW handoff
# attempt file (a6d5f17a-6127-422a-a9c7-e2c1c8202ce2.md) with a failed
# status, no handoff artifact (.xlf) produced, and the per-language handoff
# bookkeeping reset back to "nothing has shipped yet".

$wb = $excel.ActiveWorkbook

$mdFileName  = "a6d5f17a-6127-422a-a9c7-e2c1c8202ce2.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/e86378350ce4e3fb88c6da7e7533ea4af370b306/e2e/$mdFileName"
$cfgFileName = ".localization-config"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/e86378350ce4e3fb88c6da7e7533ea4af370b306/$cfgFileName"

$statusFailed = "Handoff transform failed"
$epoch        = "0001-01-01 00:00:00"

# ------------------------------------------------------------------
# Overview sheet: just the file-name hyperlink + rolled-up status
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $mdFileName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $cfgUrl, "", "", $cfgFileName)

$wsOverview.Range("B2").Value = $statusFailed
$wsOverview.Range("C2").Value = $statusFailed

# ------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de): same shape of edit on each
# ------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # The previous handoff (.xlf) artifact is gone - the transform never
    # produced one, so the "Latest Handoff File" cell/hyperlink is cleared.
    $ws.Range("C2").Clear()

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $mdFileName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", $cfgFileName)

    $ws.Range("B2").Value = $statusFailed

    # Handoff/handback dates reset to the zero-value epoch stamp.
    $ws.Range("D2").Value = $epoch
    $ws.Range("G2").Value = $epoch
    $ws.Range("D3").Value = $epoch
    $ws.Range("G3").Value = $epoch

    # Handoff reason moves from "Include" to "Ignored" for row 2; row 3 stays "Ignored".
    $ws.Range("H2").Value = "Ignored"
    $ws.Range("H3").Value = "Ignored"
}
